$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 2 banding (copy format from row 3, which already carries the
# desired "orange" band, onto row 2) -----------------------------------
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)

# --- New header cells F1:I1 --------------------------------------------
$ws.Range("F1").Value = "Manufacturer Part Number (MPN)"
$ws.Range("G1").Value = "Ordercode Farnell"
$ws.Range("H1").Value = "Ordercode Mouser"
$ws.Range("I1").Value = "Type"

# --- New data in columns F, G, H, I -------------------------------------
$ws.Range("F5").Value = 30062114
$ws.Range("H5").Value = "979-3006.2114"
$ws.Range("I5").Value = "Not specified"

$ws.Range("F6").Value = "M20-9980345"
$ws.Range("G6").Value = 1022230
$ws.Range("I6").Value = "Per stuk"

$ws.Range("F7").Value = "C0805C104M4RACAUTO"
$ws.Range("G7").Value = 3510057
$ws.Range("I7").Value = "Reel"

$ws.Range("F8").Value = "MCWR08X1001FTL"
$ws.Range("G8").Value = 2446904
$ws.Range("I8").Value = "Reel"

$ws.Range("F9").Value = "KP-2012QBC-D"
$ws.Range("G9").Value = 2217974
$ws.Range("I9").Value = "Tape en reel, afgesneden"

$ws.Range("F10").Value = "MC01W08051100K"
$ws.Range("G10").Value = 2129276
$ws.Range("I10").Value = "Reel"

$ws.Range("F12").Value = "PIC16F18446-I/SS"
$ws.Range("G12").Value = 3631504
$ws.Range("I12").Value = "Not specified"

$ws.Range("F13").Value = "MC01W0805110K"
$ws.Range("G13").Value = 2129195
$ws.Range("I13").Value = "Reel"

# --- Styling -------------------------------------------------------------
# Style used for header row F1:I1, and for the "Type" column I6:I13 + I5:
# default font, centered.
$headerType = $ws.Range("F1:I1,I5,I6,I7,I8,I9,I10,I12,I13")
$headerType.HorizontalAlignment = -4108

# Style 4 equivalent: Verdana 8pt FF333333, centered - used for the bulk
# of the Farnell/Mouser ordercode cells.
$verdanaCells = $ws.Range("F6:G13")
$verdanaCells.Font.Name = "Verdana"
$verdanaCells.Font.Size = 8
$verdanaCells.Font.Color = 3355443
$verdanaCells.HorizontalAlignment = -4108

# Row 5 uses Arial 10pt FF333333 instead (style 6 on F5 w/ thousands
# number format, style 7 on H5 without).
$arialCells = $ws.Range("F5,H5")
$arialCells.Font.Name = "Arial"
$arialCells.Font.Size = 10
$arialCells.Font.Color = 3355443
$arialCells.HorizontalAlignment = -4108

$ws.Range("F5").NumberFormat = "#,##0"

# --- Column widths --------------------------------------------------------
$ws.Range("E1").ColumnWidth = 36.85546875
$ws.Range("F1").ColumnWidth = 24.42578125
$ws.Range("G1").ColumnWidth = 21.7109375
$ws.Range("H1").ColumnWidth = 23.28515625
$ws.Range("I1").ColumnWidth = 28.42578125

# --- Selection / page setup -----------------------------------------------
$ws.Range("F3").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
